$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.086979095468127
$ws.Range("D2").Value = 0.2888123221353156

$ws.Range("C3").Value = 0.7707942849999454
$ws.Range("D3").Value = 0.4490282397069658

$ws.Range("C4").Value = 1.21983697151183
$ws.Range("D4").Value = 0.2354427166746622

$ws.Range("C5").Value = 6.371553553811489
$ws.Range("D5").Value = 0.000002065684132634971

$ws.Range("C6").Value = -0.01571425815784909
$ws.Range("D6").Value = 0.9876039936595451

$ws.Range("C7").Value = 0.4520307547817693
$ws.Range("D7").Value = 0.6556721540292556

$ws.Range("C8").Value = 3.408075197681598
$ws.Range("D8").Value = 0.002522223241182564

$ws.Range("C9").Value = 0.3635935300073379
$ws.Range("D9").Value = 0.7196322249386398

$ws.Range("C10").Value = 2.349284541777692
$ws.Range("D10").Value = 0.02819790319968396

$ws.Range("C11").Value = 2.275928573368112
$ws.Range("D11").Value = 0.03293977649519375
